$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- becomes old row 4 data (A,B,E,F,G,H,Q,R)
$ws.Range("A2").Value = 111702796
$ws.Range("B2").Value = 90687
$ws.Range("E2").Value = 5964
$ws.Range("F2").Value = "Fjällig taggsvamp s.str."
$ws.Range("G2").Value = "Sarcodon imbricatus s.str."
$ws.Range("H2").Value = "(L.:Fr.) P.Karst."
$ws.Range("Q2").Value = 516756.4585669422
$ws.Range("R2").Value = 6574760.884648616

# Row 3 <- becomes old row 2 data (A,Q,R); B,E,F,G,H unchanged (already same as old row 2)
$ws.Range("A3").Value = 111702873
$ws.Range("Q3").Value = 516761.0073171449
$ws.Range("R3").Value = 6574773.157834023

# Row 4 <- becomes old row 3 data (A,B,E,F,G,H,Q,R)
$ws.Range("A4").Value = 111702802
$ws.Range("B4").Value = 90332
$ws.Range("E4").Value = 4769
$ws.Range("F4").Value = "Svavelriska"
$ws.Range("G4").Value = "Lactarius scrobiculatus"
$ws.Range("H4").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q4").Value = 516752.3535787854
$ws.Range("R4").Value = 6574763.929792823
